# Apply "Update latest output (run 237)" changes to the workbook.
# Sheet "Schedule" (sheet1): update row 2, append new row 3.
# Sheet "Detailed" (sheet2): update rows 14-49 (B/C/E), append new rows 50-97.

$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# ---- Schedule sheet: update existing row 2 ----
$wsSchedule.Cells.Item(2,1).Value = 46067.33333333334
$wsSchedule.Cells.Item(2,2).Value = 46067.79166666666
$wsSchedule.Cells.Item(2,5).Value = 422.7723045
$wsSchedule.Cells.Item(2,6).Value = 10.16768409090909

# ---- Schedule sheet: append new row 3 (copy date/time format from row 2) ----
$wsSchedule.Cells.Item(3,1).Value = 46068.29166666666
$wsSchedule.Cells.Item(3,2).Value = 46068.75
$wsSchedule.Cells.Item(3,3).Value = 11
$wsSchedule.Cells.Item(3,4).Value = 41.58
$wsSchedule.Cells.Item(3,5).Value = 776.3698409999997
$wsSchedule.Cells.Item(3,6).Value = 18.67171334776334
$wsSchedule.Cells.Item(3,1).NumberFormat = $wsSchedule.Cells.Item(2,1).NumberFormat
$wsSchedule.Cells.Item(3,2).NumberFormat = $wsSchedule.Cells.Item(2,2).NumberFormat

# ---- Detailed sheet: update existing rows 14-49 (Price / Type / Pump_Status) ----
$wsDetailed.Cells.Item(14,2).Value = 84.79
$wsDetailed.Cells.Item(15,2).Value = 84.79
$wsDetailed.Cells.Item(16,2).Value = 73.19
$wsDetailed.Cells.Item(16,3).Value = "historical"
$wsDetailed.Cells.Item(17,2).Value = 51.59311
$wsDetailed.Cells.Item(17,3).Value = "historical"
$wsDetailed.Cells.Item(17,5).Value = "OFF"
$wsDetailed.Cells.Item(18,2).Value = 56.98
$wsDetailed.Cells.Item(18,3).Value = "historical"
$wsDetailed.Cells.Item(19,2).Value = 47.45981
$wsDetailed.Cells.Item(19,3).Value = "historical"
$wsDetailed.Cells.Item(20,2).Value = 36.07
$wsDetailed.Cells.Item(20,3).Value = "historical"
$wsDetailed.Cells.Item(21,2).Value = 31.1147
$wsDetailed.Cells.Item(21,3).Value = "historical"
$wsDetailed.Cells.Item(22,2).Value = 22.07
$wsDetailed.Cells.Item(22,3).Value = "historical"
$wsDetailed.Cells.Item(23,2).Value = 30.94737
$wsDetailed.Cells.Item(23,3).Value = "historical"
$wsDetailed.Cells.Item(24,2).Value = 30.52819
$wsDetailed.Cells.Item(24,3).Value = "historical"
$wsDetailed.Cells.Item(25,2).Value = 22.07
$wsDetailed.Cells.Item(25,3).Value = "historical"
$wsDetailed.Cells.Item(26,2).Value = 8.03814
$wsDetailed.Cells.Item(26,3).Value = "historical"
$wsDetailed.Cells.Item(27,2).Value = -4.68971
$wsDetailed.Cells.Item(27,3).Value = "historical"
$wsDetailed.Cells.Item(28,2).Value = -0.99369
$wsDetailed.Cells.Item(28,3).Value = "historical"
$wsDetailed.Cells.Item(29,2).Value = -2.66469
$wsDetailed.Cells.Item(29,3).Value = "historical"
$wsDetailed.Cells.Item(30,2).Value = 2.9527
$wsDetailed.Cells.Item(30,3).Value = "historical"
$wsDetailed.Cells.Item(31,2).Value = 0.51
$wsDetailed.Cells.Item(31,3).Value = "historical"
$wsDetailed.Cells.Item(32,2).Value = -3.75985
$wsDetailed.Cells.Item(32,3).Value = "historical"
$wsDetailed.Cells.Item(33,2).Value = -5.01
$wsDetailed.Cells.Item(34,2).Value = 0.01723
$wsDetailed.Cells.Item(35,2).Value = -3.75985
$wsDetailed.Cells.Item(36,2).Value = 36.0601
$wsDetailed.Cells.Item(37,2).Value = 36.0601
$wsDetailed.Cells.Item(38,2).Value = 45.11991
$wsDetailed.Cells.Item(39,2).Value = 48.49216
$wsDetailed.Cells.Item(39,5).Value = "ON"
$wsDetailed.Cells.Item(40,2).Value = 57.17169
$wsDetailed.Cells.Item(45,2).Value = 57.31
$wsDetailed.Cells.Item(48,2).Value = 56.98
$wsDetailed.Cells.Item(49,2).Value = 56.98

# ---- Detailed sheet: append new rows 50-97 for 2026-02-15 (date 46068) ----
$wsDetailed.Cells.Item(50,1).Value = 46068
$wsDetailed.Cells.Item(50,2).Value = 52.21838
$wsDetailed.Cells.Item(50,3).Value = "forecast"
$wsDetailed.Cells.Item(50,4).Value = 46068
$wsDetailed.Cells.Item(50,5).Value = "OFF"
$wsDetailed.Cells.Item(50,1).NumberFormat = $wsDetailed.Cells.Item(49,1).NumberFormat
$wsDetailed.Cells.Item(50,4).NumberFormat = $wsDetailed.Cells.Item(49,4).NumberFormat
$wsDetailed.Cells.Item(51,1).Value = 46068.02083333334
$wsDetailed.Cells.Item(51,2).Value = 56.98
$wsDetailed.Cells.Item(51,3).Value = "forecast"
$wsDetailed.Cells.Item(51,4).Value = 46068
$wsDetailed.Cells.Item(51,5).Value = "OFF"
$wsDetailed.Cells.Item(51,1).NumberFormat = $wsDetailed.Cells.Item(50,1).NumberFormat
$wsDetailed.Cells.Item(51,4).NumberFormat = $wsDetailed.Cells.Item(50,4).NumberFormat
$wsDetailed.Cells.Item(52,1).Value = 46068.04166666666
$wsDetailed.Cells.Item(52,2).Value = 56.98
$wsDetailed.Cells.Item(52,3).Value = "forecast"
$wsDetailed.Cells.Item(52,4).Value = 46068
$wsDetailed.Cells.Item(52,5).Value = "OFF"
$wsDetailed.Cells.Item(52,1).NumberFormat = $wsDetailed.Cells.Item(51,1).NumberFormat
$wsDetailed.Cells.Item(52,4).NumberFormat = $wsDetailed.Cells.Item(51,4).NumberFormat
$wsDetailed.Cells.Item(53,1).Value = 46068.0625
$wsDetailed.Cells.Item(53,2).Value = 56.98
$wsDetailed.Cells.Item(53,3).Value = "forecast"
$wsDetailed.Cells.Item(53,4).Value = 46068
$wsDetailed.Cells.Item(53,5).Value = "OFF"
$wsDetailed.Cells.Item(53,1).NumberFormat = $wsDetailed.Cells.Item(52,1).NumberFormat
$wsDetailed.Cells.Item(53,4).NumberFormat = $wsDetailed.Cells.Item(52,4).NumberFormat
$wsDetailed.Cells.Item(54,1).Value = 46068.08333333334
$wsDetailed.Cells.Item(54,2).Value = 54.33483
$wsDetailed.Cells.Item(54,3).Value = "forecast"
$wsDetailed.Cells.Item(54,4).Value = 46068
$wsDetailed.Cells.Item(54,5).Value = "OFF"
$wsDetailed.Cells.Item(54,1).NumberFormat = $wsDetailed.Cells.Item(53,1).NumberFormat
$wsDetailed.Cells.Item(54,4).NumberFormat = $wsDetailed.Cells.Item(53,4).NumberFormat
$wsDetailed.Cells.Item(55,1).Value = 46068.10416666666
$wsDetailed.Cells.Item(55,2).Value = 54.72429
$wsDetailed.Cells.Item(55,3).Value = "forecast"
$wsDetailed.Cells.Item(55,4).Value = 46068
$wsDetailed.Cells.Item(55,5).Value = "OFF"
$wsDetailed.Cells.Item(55,1).NumberFormat = $wsDetailed.Cells.Item(54,1).NumberFormat
$wsDetailed.Cells.Item(55,4).NumberFormat = $wsDetailed.Cells.Item(54,4).NumberFormat
$wsDetailed.Cells.Item(56,1).Value = 46068.125
$wsDetailed.Cells.Item(56,2).Value = 54.76507
$wsDetailed.Cells.Item(56,3).Value = "forecast"
$wsDetailed.Cells.Item(56,4).Value = 46068
$wsDetailed.Cells.Item(56,5).Value = "OFF"
$wsDetailed.Cells.Item(56,1).NumberFormat = $wsDetailed.Cells.Item(55,1).NumberFormat
$wsDetailed.Cells.Item(56,4).NumberFormat = $wsDetailed.Cells.Item(55,4).NumberFormat
$wsDetailed.Cells.Item(57,1).Value = 46068.14583333334
$wsDetailed.Cells.Item(57,2).Value = 56.83994
$wsDetailed.Cells.Item(57,3).Value = "forecast"
$wsDetailed.Cells.Item(57,4).Value = 46068
$wsDetailed.Cells.Item(57,5).Value = "OFF"
$wsDetailed.Cells.Item(57,1).NumberFormat = $wsDetailed.Cells.Item(56,1).NumberFormat
$wsDetailed.Cells.Item(57,4).NumberFormat = $wsDetailed.Cells.Item(56,4).NumberFormat
$wsDetailed.Cells.Item(58,1).Value = 46068.16666666666
$wsDetailed.Cells.Item(58,2).Value = 56.98
$wsDetailed.Cells.Item(58,3).Value = "forecast"
$wsDetailed.Cells.Item(58,4).Value = 46068
$wsDetailed.Cells.Item(58,5).Value = "OFF"
$wsDetailed.Cells.Item(58,1).NumberFormat = $wsDetailed.Cells.Item(57,1).NumberFormat
$wsDetailed.Cells.Item(58,4).NumberFormat = $wsDetailed.Cells.Item(57,4).NumberFormat
$wsDetailed.Cells.Item(59,1).Value = 46068.1875
$wsDetailed.Cells.Item(59,2).Value = 50.80126
$wsDetailed.Cells.Item(59,3).Value = "forecast"
$wsDetailed.Cells.Item(59,4).Value = 46068
$wsDetailed.Cells.Item(59,5).Value = "OFF"
$wsDetailed.Cells.Item(59,1).NumberFormat = $wsDetailed.Cells.Item(58,1).NumberFormat
$wsDetailed.Cells.Item(59,4).NumberFormat = $wsDetailed.Cells.Item(58,4).NumberFormat
$wsDetailed.Cells.Item(60,1).Value = 46068.20833333334
$wsDetailed.Cells.Item(60,2).Value = 51.08547
$wsDetailed.Cells.Item(60,3).Value = "forecast"
$wsDetailed.Cells.Item(60,4).Value = 46068
$wsDetailed.Cells.Item(60,5).Value = "OFF"
$wsDetailed.Cells.Item(60,1).NumberFormat = $wsDetailed.Cells.Item(59,1).NumberFormat
$wsDetailed.Cells.Item(60,4).NumberFormat = $wsDetailed.Cells.Item(59,4).NumberFormat
$wsDetailed.Cells.Item(61,1).Value = 46068.22916666666
$wsDetailed.Cells.Item(61,2).Value = 54.69654
$wsDetailed.Cells.Item(61,3).Value = "forecast"
$wsDetailed.Cells.Item(61,4).Value = 46068
$wsDetailed.Cells.Item(61,5).Value = "OFF"
$wsDetailed.Cells.Item(61,1).NumberFormat = $wsDetailed.Cells.Item(60,1).NumberFormat
$wsDetailed.Cells.Item(61,4).NumberFormat = $wsDetailed.Cells.Item(60,4).NumberFormat
$wsDetailed.Cells.Item(62,1).Value = 46068.25
$wsDetailed.Cells.Item(62,2).Value = 52.98916
$wsDetailed.Cells.Item(62,3).Value = "forecast"
$wsDetailed.Cells.Item(62,4).Value = 46068
$wsDetailed.Cells.Item(62,5).Value = "OFF"
$wsDetailed.Cells.Item(62,1).NumberFormat = $wsDetailed.Cells.Item(61,1).NumberFormat
$wsDetailed.Cells.Item(62,4).NumberFormat = $wsDetailed.Cells.Item(61,4).NumberFormat
$wsDetailed.Cells.Item(63,1).Value = 46068.27083333334
$wsDetailed.Cells.Item(63,2).Value = 56.80005
$wsDetailed.Cells.Item(63,3).Value = "forecast"
$wsDetailed.Cells.Item(63,4).Value = 46068
$wsDetailed.Cells.Item(63,5).Value = "OFF"
$wsDetailed.Cells.Item(63,1).NumberFormat = $wsDetailed.Cells.Item(62,1).NumberFormat
$wsDetailed.Cells.Item(63,4).NumberFormat = $wsDetailed.Cells.Item(62,4).NumberFormat
$wsDetailed.Cells.Item(64,1).Value = 46068.29166666666
$wsDetailed.Cells.Item(64,2).Value = 36.0601
$wsDetailed.Cells.Item(64,3).Value = "forecast"
$wsDetailed.Cells.Item(64,4).Value = 46068
$wsDetailed.Cells.Item(64,5).Value = "ON"
$wsDetailed.Cells.Item(64,1).NumberFormat = $wsDetailed.Cells.Item(63,1).NumberFormat
$wsDetailed.Cells.Item(64,4).NumberFormat = $wsDetailed.Cells.Item(63,4).NumberFormat
$wsDetailed.Cells.Item(65,1).Value = 46068.3125
$wsDetailed.Cells.Item(65,2).Value = 36.0601
$wsDetailed.Cells.Item(65,3).Value = "forecast"
$wsDetailed.Cells.Item(65,4).Value = 46068
$wsDetailed.Cells.Item(65,5).Value = "ON"
$wsDetailed.Cells.Item(65,1).NumberFormat = $wsDetailed.Cells.Item(64,1).NumberFormat
$wsDetailed.Cells.Item(65,4).NumberFormat = $wsDetailed.Cells.Item(64,4).NumberFormat
$wsDetailed.Cells.Item(66,1).Value = 46068.33333333334
$wsDetailed.Cells.Item(66,2).Value = 36.0601
$wsDetailed.Cells.Item(66,3).Value = "forecast"
$wsDetailed.Cells.Item(66,4).Value = 46068
$wsDetailed.Cells.Item(66,5).Value = "ON"
$wsDetailed.Cells.Item(66,1).NumberFormat = $wsDetailed.Cells.Item(65,1).NumberFormat
$wsDetailed.Cells.Item(66,4).NumberFormat = $wsDetailed.Cells.Item(65,4).NumberFormat
$wsDetailed.Cells.Item(67,1).Value = 46068.35416666666
$wsDetailed.Cells.Item(67,2).Value = 36.0601
$wsDetailed.Cells.Item(67,3).Value = "forecast"
$wsDetailed.Cells.Item(67,4).Value = 46068
$wsDetailed.Cells.Item(67,5).Value = "ON"
$wsDetailed.Cells.Item(67,1).NumberFormat = $wsDetailed.Cells.Item(66,1).NumberFormat
$wsDetailed.Cells.Item(67,4).NumberFormat = $wsDetailed.Cells.Item(66,4).NumberFormat
$wsDetailed.Cells.Item(68,1).Value = 46068.375
$wsDetailed.Cells.Item(68,2).Value = 36.07
$wsDetailed.Cells.Item(68,3).Value = "forecast"
$wsDetailed.Cells.Item(68,4).Value = 46068
$wsDetailed.Cells.Item(68,5).Value = "ON"
$wsDetailed.Cells.Item(68,1).NumberFormat = $wsDetailed.Cells.Item(67,1).NumberFormat
$wsDetailed.Cells.Item(68,4).NumberFormat = $wsDetailed.Cells.Item(67,4).NumberFormat
$wsDetailed.Cells.Item(69,1).Value = 46068.39583333334
$wsDetailed.Cells.Item(69,2).Value = 36.0601
$wsDetailed.Cells.Item(69,3).Value = "forecast"
$wsDetailed.Cells.Item(69,4).Value = 46068
$wsDetailed.Cells.Item(69,5).Value = "ON"
$wsDetailed.Cells.Item(69,1).NumberFormat = $wsDetailed.Cells.Item(68,1).NumberFormat
$wsDetailed.Cells.Item(69,4).NumberFormat = $wsDetailed.Cells.Item(68,4).NumberFormat
$wsDetailed.Cells.Item(70,1).Value = 46068.41666666666
$wsDetailed.Cells.Item(70,2).Value = 36.0601
$wsDetailed.Cells.Item(70,3).Value = "forecast"
$wsDetailed.Cells.Item(70,4).Value = 46068
$wsDetailed.Cells.Item(70,5).Value = "ON"
$wsDetailed.Cells.Item(70,1).NumberFormat = $wsDetailed.Cells.Item(69,1).NumberFormat
$wsDetailed.Cells.Item(70,4).NumberFormat = $wsDetailed.Cells.Item(69,4).NumberFormat
$wsDetailed.Cells.Item(71,1).Value = 46068.4375
$wsDetailed.Cells.Item(71,2).Value = 36.07
$wsDetailed.Cells.Item(71,3).Value = "forecast"
$wsDetailed.Cells.Item(71,4).Value = 46068
$wsDetailed.Cells.Item(71,5).Value = "ON"
$wsDetailed.Cells.Item(71,1).NumberFormat = $wsDetailed.Cells.Item(70,1).NumberFormat
$wsDetailed.Cells.Item(71,4).NumberFormat = $wsDetailed.Cells.Item(70,4).NumberFormat
$wsDetailed.Cells.Item(72,1).Value = 46068.45833333334
$wsDetailed.Cells.Item(72,2).Value = 36.0601
$wsDetailed.Cells.Item(72,3).Value = "forecast"
$wsDetailed.Cells.Item(72,4).Value = 46068
$wsDetailed.Cells.Item(72,5).Value = "ON"
$wsDetailed.Cells.Item(72,1).NumberFormat = $wsDetailed.Cells.Item(71,1).NumberFormat
$wsDetailed.Cells.Item(72,4).NumberFormat = $wsDetailed.Cells.Item(71,4).NumberFormat
$wsDetailed.Cells.Item(73,1).Value = 46068.47916666666
$wsDetailed.Cells.Item(73,2).Value = 36.0601
$wsDetailed.Cells.Item(73,3).Value = "forecast"
$wsDetailed.Cells.Item(73,4).Value = 46068
$wsDetailed.Cells.Item(73,5).Value = "ON"
$wsDetailed.Cells.Item(73,1).NumberFormat = $wsDetailed.Cells.Item(72,1).NumberFormat
$wsDetailed.Cells.Item(73,4).NumberFormat = $wsDetailed.Cells.Item(72,4).NumberFormat
$wsDetailed.Cells.Item(74,1).Value = 46068.5
$wsDetailed.Cells.Item(74,2).Value = 36.0601
$wsDetailed.Cells.Item(74,3).Value = "forecast"
$wsDetailed.Cells.Item(74,4).Value = 46068
$wsDetailed.Cells.Item(74,5).Value = "ON"
$wsDetailed.Cells.Item(74,1).NumberFormat = $wsDetailed.Cells.Item(73,1).NumberFormat
$wsDetailed.Cells.Item(74,4).NumberFormat = $wsDetailed.Cells.Item(73,4).NumberFormat
$wsDetailed.Cells.Item(75,1).Value = 46068.52083333334
$wsDetailed.Cells.Item(75,2).Value = 28.91993
$wsDetailed.Cells.Item(75,3).Value = "forecast"
$wsDetailed.Cells.Item(75,4).Value = 46068
$wsDetailed.Cells.Item(75,5).Value = "ON"
$wsDetailed.Cells.Item(75,1).NumberFormat = $wsDetailed.Cells.Item(74,1).NumberFormat
$wsDetailed.Cells.Item(75,4).NumberFormat = $wsDetailed.Cells.Item(74,4).NumberFormat
$wsDetailed.Cells.Item(76,1).Value = 46068.54166666666
$wsDetailed.Cells.Item(76,2).Value = 34.97618
$wsDetailed.Cells.Item(76,3).Value = "forecast"
$wsDetailed.Cells.Item(76,4).Value = 46068
$wsDetailed.Cells.Item(76,5).Value = "ON"
$wsDetailed.Cells.Item(76,1).NumberFormat = $wsDetailed.Cells.Item(75,1).NumberFormat
$wsDetailed.Cells.Item(76,4).NumberFormat = $wsDetailed.Cells.Item(75,4).NumberFormat
$wsDetailed.Cells.Item(77,1).Value = 46068.5625
$wsDetailed.Cells.Item(77,2).Value = 33.20999
$wsDetailed.Cells.Item(77,3).Value = "forecast"
$wsDetailed.Cells.Item(77,4).Value = 46068
$wsDetailed.Cells.Item(77,5).Value = "ON"
$wsDetailed.Cells.Item(77,1).NumberFormat = $wsDetailed.Cells.Item(76,1).NumberFormat
$wsDetailed.Cells.Item(77,4).NumberFormat = $wsDetailed.Cells.Item(76,4).NumberFormat
$wsDetailed.Cells.Item(78,1).Value = 46068.58333333334
$wsDetailed.Cells.Item(78,2).Value = 36.0601
$wsDetailed.Cells.Item(78,3).Value = "forecast"
$wsDetailed.Cells.Item(78,4).Value = 46068
$wsDetailed.Cells.Item(78,5).Value = "ON"
$wsDetailed.Cells.Item(78,1).NumberFormat = $wsDetailed.Cells.Item(77,1).NumberFormat
$wsDetailed.Cells.Item(78,4).NumberFormat = $wsDetailed.Cells.Item(77,4).NumberFormat
$wsDetailed.Cells.Item(79,1).Value = 46068.60416666666
$wsDetailed.Cells.Item(79,2).Value = 36.0601
$wsDetailed.Cells.Item(79,3).Value = "forecast"
$wsDetailed.Cells.Item(79,4).Value = 46068
$wsDetailed.Cells.Item(79,5).Value = "ON"
$wsDetailed.Cells.Item(79,1).NumberFormat = $wsDetailed.Cells.Item(78,1).NumberFormat
$wsDetailed.Cells.Item(79,4).NumberFormat = $wsDetailed.Cells.Item(78,4).NumberFormat
$wsDetailed.Cells.Item(80,1).Value = 46068.625
$wsDetailed.Cells.Item(80,2).Value = 45.52957
$wsDetailed.Cells.Item(80,3).Value = "forecast"
$wsDetailed.Cells.Item(80,4).Value = 46068
$wsDetailed.Cells.Item(80,5).Value = "ON"
$wsDetailed.Cells.Item(80,1).NumberFormat = $wsDetailed.Cells.Item(79,1).NumberFormat
$wsDetailed.Cells.Item(80,4).NumberFormat = $wsDetailed.Cells.Item(79,4).NumberFormat
$wsDetailed.Cells.Item(81,1).Value = 46068.64583333334
$wsDetailed.Cells.Item(81,2).Value = 48.03047
$wsDetailed.Cells.Item(81,3).Value = "forecast"
$wsDetailed.Cells.Item(81,4).Value = 46068
$wsDetailed.Cells.Item(81,5).Value = "ON"
$wsDetailed.Cells.Item(81,1).NumberFormat = $wsDetailed.Cells.Item(80,1).NumberFormat
$wsDetailed.Cells.Item(81,4).NumberFormat = $wsDetailed.Cells.Item(80,4).NumberFormat
$wsDetailed.Cells.Item(82,1).Value = 46068.66666666666
$wsDetailed.Cells.Item(82,2).Value = 45.5218
$wsDetailed.Cells.Item(82,3).Value = "forecast"
$wsDetailed.Cells.Item(82,4).Value = 46068
$wsDetailed.Cells.Item(82,5).Value = "ON"
$wsDetailed.Cells.Item(82,1).NumberFormat = $wsDetailed.Cells.Item(81,1).NumberFormat
$wsDetailed.Cells.Item(82,4).NumberFormat = $wsDetailed.Cells.Item(81,4).NumberFormat
$wsDetailed.Cells.Item(83,1).Value = 46068.6875
$wsDetailed.Cells.Item(83,2).Value = 29.924
$wsDetailed.Cells.Item(83,3).Value = "forecast"
$wsDetailed.Cells.Item(83,4).Value = 46068
$wsDetailed.Cells.Item(83,5).Value = "ON"
$wsDetailed.Cells.Item(83,1).NumberFormat = $wsDetailed.Cells.Item(82,1).NumberFormat
$wsDetailed.Cells.Item(83,4).NumberFormat = $wsDetailed.Cells.Item(82,4).NumberFormat
$wsDetailed.Cells.Item(84,1).Value = 46068.70833333334
$wsDetailed.Cells.Item(84,2).Value = 29.94701
$wsDetailed.Cells.Item(84,3).Value = "forecast"
$wsDetailed.Cells.Item(84,4).Value = 46068
$wsDetailed.Cells.Item(84,5).Value = "ON"
$wsDetailed.Cells.Item(84,1).NumberFormat = $wsDetailed.Cells.Item(83,1).NumberFormat
$wsDetailed.Cells.Item(84,4).NumberFormat = $wsDetailed.Cells.Item(83,4).NumberFormat
$wsDetailed.Cells.Item(85,1).Value = 46068.72916666666
$wsDetailed.Cells.Item(85,2).Value = 31.41671
$wsDetailed.Cells.Item(85,3).Value = "forecast"
$wsDetailed.Cells.Item(85,4).Value = 46068
$wsDetailed.Cells.Item(85,5).Value = "ON"
$wsDetailed.Cells.Item(85,1).NumberFormat = $wsDetailed.Cells.Item(84,1).NumberFormat
$wsDetailed.Cells.Item(85,4).NumberFormat = $wsDetailed.Cells.Item(84,4).NumberFormat
$wsDetailed.Cells.Item(86,1).Value = 46068.75
$wsDetailed.Cells.Item(86,2).Value = 46.9392
$wsDetailed.Cells.Item(86,3).Value = "forecast"
$wsDetailed.Cells.Item(86,4).Value = 46068
$wsDetailed.Cells.Item(86,5).Value = "OFF"
$wsDetailed.Cells.Item(86,1).NumberFormat = $wsDetailed.Cells.Item(85,1).NumberFormat
$wsDetailed.Cells.Item(86,4).NumberFormat = $wsDetailed.Cells.Item(85,4).NumberFormat
$wsDetailed.Cells.Item(87,1).Value = 46068.77083333334
$wsDetailed.Cells.Item(87,2).Value = 57.31
$wsDetailed.Cells.Item(87,3).Value = "forecast"
$wsDetailed.Cells.Item(87,4).Value = 46068
$wsDetailed.Cells.Item(87,5).Value = "OFF"
$wsDetailed.Cells.Item(87,1).NumberFormat = $wsDetailed.Cells.Item(86,1).NumberFormat
$wsDetailed.Cells.Item(87,4).NumberFormat = $wsDetailed.Cells.Item(86,4).NumberFormat
$wsDetailed.Cells.Item(88,1).Value = 46068.79166666666
$wsDetailed.Cells.Item(88,2).Value = 60.54926
$wsDetailed.Cells.Item(88,3).Value = "forecast"
$wsDetailed.Cells.Item(88,4).Value = 46068
$wsDetailed.Cells.Item(88,5).Value = "OFF"
$wsDetailed.Cells.Item(88,1).NumberFormat = $wsDetailed.Cells.Item(87,1).NumberFormat
$wsDetailed.Cells.Item(88,4).NumberFormat = $wsDetailed.Cells.Item(87,4).NumberFormat
$wsDetailed.Cells.Item(89,1).Value = 46068.8125
$wsDetailed.Cells.Item(89,2).Value = 64.89
$wsDetailed.Cells.Item(89,3).Value = "forecast"
$wsDetailed.Cells.Item(89,4).Value = 46068
$wsDetailed.Cells.Item(89,5).Value = "OFF"
$wsDetailed.Cells.Item(89,1).NumberFormat = $wsDetailed.Cells.Item(88,1).NumberFormat
$wsDetailed.Cells.Item(89,4).NumberFormat = $wsDetailed.Cells.Item(88,4).NumberFormat
$wsDetailed.Cells.Item(90,1).Value = 46068.83333333334
$wsDetailed.Cells.Item(90,2).Value = 64.89
$wsDetailed.Cells.Item(90,3).Value = "forecast"
$wsDetailed.Cells.Item(90,4).Value = 46068
$wsDetailed.Cells.Item(90,5).Value = "OFF"
$wsDetailed.Cells.Item(90,1).NumberFormat = $wsDetailed.Cells.Item(89,1).NumberFormat
$wsDetailed.Cells.Item(90,4).NumberFormat = $wsDetailed.Cells.Item(89,4).NumberFormat
$wsDetailed.Cells.Item(91,1).Value = 46068.85416666666
$wsDetailed.Cells.Item(91,2).Value = 60.06944
$wsDetailed.Cells.Item(91,3).Value = "forecast"
$wsDetailed.Cells.Item(91,4).Value = 46068
$wsDetailed.Cells.Item(91,5).Value = "OFF"
$wsDetailed.Cells.Item(91,1).NumberFormat = $wsDetailed.Cells.Item(90,1).NumberFormat
$wsDetailed.Cells.Item(91,4).NumberFormat = $wsDetailed.Cells.Item(90,4).NumberFormat
$wsDetailed.Cells.Item(92,1).Value = 46068.875
$wsDetailed.Cells.Item(92,2).Value = 60.14357
$wsDetailed.Cells.Item(92,3).Value = "forecast"
$wsDetailed.Cells.Item(92,4).Value = 46068
$wsDetailed.Cells.Item(92,5).Value = "OFF"
$wsDetailed.Cells.Item(92,1).NumberFormat = $wsDetailed.Cells.Item(91,1).NumberFormat
$wsDetailed.Cells.Item(92,4).NumberFormat = $wsDetailed.Cells.Item(91,4).NumberFormat
$wsDetailed.Cells.Item(93,1).Value = 46068.89583333334
$wsDetailed.Cells.Item(93,2).Value = 58.47004
$wsDetailed.Cells.Item(93,3).Value = "forecast"
$wsDetailed.Cells.Item(93,4).Value = 46068
$wsDetailed.Cells.Item(93,5).Value = "OFF"
$wsDetailed.Cells.Item(93,1).NumberFormat = $wsDetailed.Cells.Item(92,1).NumberFormat
$wsDetailed.Cells.Item(93,4).NumberFormat = $wsDetailed.Cells.Item(92,4).NumberFormat
$wsDetailed.Cells.Item(94,1).Value = 46068.91666666666
$wsDetailed.Cells.Item(94,2).Value = 57.09
$wsDetailed.Cells.Item(94,3).Value = "forecast"
$wsDetailed.Cells.Item(94,4).Value = 46068
$wsDetailed.Cells.Item(94,5).Value = "OFF"
$wsDetailed.Cells.Item(94,1).NumberFormat = $wsDetailed.Cells.Item(93,1).NumberFormat
$wsDetailed.Cells.Item(94,4).NumberFormat = $wsDetailed.Cells.Item(93,4).NumberFormat
$wsDetailed.Cells.Item(95,1).Value = 46068.9375
$wsDetailed.Cells.Item(95,2).Value = 48.97061
$wsDetailed.Cells.Item(95,3).Value = "forecast"
$wsDetailed.Cells.Item(95,4).Value = 46068
$wsDetailed.Cells.Item(95,5).Value = "OFF"
$wsDetailed.Cells.Item(95,1).NumberFormat = $wsDetailed.Cells.Item(94,1).NumberFormat
$wsDetailed.Cells.Item(95,4).NumberFormat = $wsDetailed.Cells.Item(94,4).NumberFormat
$wsDetailed.Cells.Item(96,1).Value = 46068.95833333334
$wsDetailed.Cells.Item(96,2).Value = 45.92645
$wsDetailed.Cells.Item(96,3).Value = "forecast"
$wsDetailed.Cells.Item(96,4).Value = 46068
$wsDetailed.Cells.Item(96,5).Value = "OFF"
$wsDetailed.Cells.Item(96,1).NumberFormat = $wsDetailed.Cells.Item(95,1).NumberFormat
$wsDetailed.Cells.Item(96,4).NumberFormat = $wsDetailed.Cells.Item(95,4).NumberFormat
$wsDetailed.Cells.Item(97,1).Value = 46068.97916666666
$wsDetailed.Cells.Item(97,2).Value = 36.0601
$wsDetailed.Cells.Item(97,3).Value = "forecast"
$wsDetailed.Cells.Item(97,4).Value = 46068
$wsDetailed.Cells.Item(97,5).Value = "OFF"
$wsDetailed.Cells.Item(97,1).NumberFormat = $wsDetailed.Cells.Item(96,1).NumberFormat
$wsDetailed.Cells.Item(97,4).NumberFormat = $wsDetailed.Cells.Item(96,4).NumberFormat

Write-Host "Applied run 237 update."
